# Insert a new worksheet "Digestion1" right after "Samples", matching the
# commit that added a digestion-yield summary tab (sheetId 4) ahead of the
# existing "Well layout1"/"Well layout2" tabs.
$wb = $excel.ActiveWorkbook
$samples = $wb.Worksheets.Item("Samples")
$new = $wb.Worksheets.Add($null, $samples)
$new.Name = "Digestion1"

# Row 1 headers 1-12
$new.Range("B1").Value = 1
$new.Range("C1").Value = 2
$new.Range("D1").Value = 3
$new.Range("E1").Value = 4
$new.Range("F1").Value = 5
$new.Range("G1").Value = 6
$new.Range("H1").Value = 7
$new.Range("I1").Value = 8
$new.Range("J1").Value = 9
$new.Range("K1").Value = 10
$new.Range("L1").Value = 11
$new.Range("M1").Value = 12

# Row labels A-H in column A, rows 2-9
$new.Range("A2").Value = "A"
$new.Range("A3").Value = "B"
$new.Range("A4").Value = "C"
$new.Range("A5").Value = "D"
$new.Range("A6").Value = "E"
$new.Range("A7").Value = "F"
$new.Range("A8").Value = "G"
$new.Range("A9").Value = "H"

# Data grid rows 3-8, columns C-L
$new.Range("C3").Value = "1,HP2,SZ-0526-1-1,sample,526,1,0526-1,1,1653.415,1641.155,12.26,SZ,,HPL,V2P1,07-27-2023"
$new.Range("D3").Value = "1,HP2,SZ-0526-1-1,sample,526,1,0526-1,1,1653.415,1641.155,12.26,SZ,,HPL,V2P1,07-27-2023"
$new.Range("E3").Value = "1,HP2,SZ-0526-1-1,sample,526,1,0526-1,1,1653.415,1641.155,12.26,SZ,,HPL,V2P1,07-27-2023"
$new.Range("F3").Value = "1,HP2,SZ-0526-1-1,sample,526,1,0526-1,1,1653.415,1641.155,12.26,SZ,,HPL,V2P1,07-27-2023"
$new.Range("G3").Value = "7,HP2,SZ-0526-3-1,sample,526,3,0526-3,1,1640.535,1627.705,12.83,SZ,,FBS,V2P1,07-26-2023"
$new.Range("H3").Value = "7,HP2,SZ-0526-3-1,sample,526,3,0526-3,1,1640.535,1627.705,12.83,SZ,,FBS,V2P1,07-26-2023"
$new.Range("I3").Value = "7,HP2,SZ-0526-3-1,sample,526,3,0526-3,1,1640.535,1627.705,12.83,SZ,,FBS,V2P1,07-26-2023"
$new.Range("J3").Value = "7,HP2,SZ-0526-3-1,sample,526,3,0526-3,1,1640.535,1627.705,12.83,SZ,,FBS,V2P1,07-26-2023"
$new.Range("K3").Value = "13,HP2,SZ-0602-1-1,sample,602,1,0602-1,1,1646.415,1634.67,11.745,SZ,,HPL,PLA,07-26-2023"
$new.Range("L3").Value = "13,HP2,SZ-0602-1-1,sample,602,1,0602-1,1,1646.415,1634.67,11.745,SZ,,HPL,PLA,07-26-2023"
$new.Range("C4").Value = "2,HP2,SZ-0526-1-2,sample,526,1,0526-1,2,1651.245,1638.63,12.615,SZ,,HPL,V2P1,07-26-2023"
$new.Range("D4").Value = "2,HP2,SZ-0526-1-2,sample,526,1,0526-1,2,1651.245,1638.63,12.615,SZ,,HPL,V2P1,07-26-2023"
$new.Range("E4").Value = "2,HP2,SZ-0526-1-2,sample,526,1,0526-1,2,1651.245,1638.63,12.615,SZ,,HPL,V2P1,07-26-2023"
$new.Range("F4").Value = "2,HP2,SZ-0526-1-2,sample,526,1,0526-1,2,1651.245,1638.63,12.615,SZ,,HPL,V2P1,07-26-2023"
$new.Range("G4").Value = "8,HP2,SZ-0526-3-2,sample,526,3,0526-3,2,1654.8,1642.14,12.66,SZ,,FBS,V2P1,07-26-2023"
$new.Range("H4").Value = "8,HP2,SZ-0526-3-2,sample,526,3,0526-3,2,1654.8,1642.14,12.66,SZ,,FBS,V2P1,07-26-2023"
$new.Range("I4").Value = "8,HP2,SZ-0526-3-2,sample,526,3,0526-3,2,1654.8,1642.14,12.66,SZ,,FBS,V2P1,07-26-2023"
$new.Range("J4").Value = "8,HP2,SZ-0526-3-2,sample,526,3,0526-3,2,1654.8,1642.14,12.66,SZ,,FBS,V2P1,07-26-2023"
$new.Range("K4").Value = "14,HP2,SZ-0602-1-2,sample,602,1,0602-1,2,1675.795,1664.3,11.495,SZ,,HPL,PLA,07-26-2023"
$new.Range("L4").Value = "14,HP2,SZ-0602-1-2,sample,602,1,0602-1,2,1675.795,1664.3,11.495,SZ,,HPL,PLA,07-26-2023"
$new.Range("C5").Value = "3,HP2,SZ-0526-1-3,sample,526,1,0526-1,3,1654.31,1641.285,13.025,SZ,,HPL,V2P1,07-26-2023"
$new.Range("D5").Value = "3,HP2,SZ-0526-1-3,sample,526,1,0526-1,3,1654.31,1641.285,13.025,SZ,,HPL,V2P1,07-26-2023"
$new.Range("E5").Value = "3,HP2,SZ-0526-1-3,sample,526,1,0526-1,3,1654.31,1641.285,13.025,SZ,,HPL,V2P1,07-26-2023"
$new.Range("F5").Value = "3,HP2,SZ-0526-1-3,sample,526,1,0526-1,3,1654.31,1641.285,13.025,SZ,,HPL,V2P1,07-26-2023"
$new.Range("G5").Value = "9,HP2,SZ-0526-3-3,sample,526,3,0526-3,3,1646.195,1632.825,13.37,SZ,,FBS,V2P1,07-26-2023"
$new.Range("H5").Value = "9,HP2,SZ-0526-3-3,sample,526,3,0526-3,3,1646.195,1632.825,13.37,SZ,,FBS,V2P1,07-26-2023"
$new.Range("I5").Value = "9,HP2,SZ-0526-3-3,sample,526,3,0526-3,3,1646.195,1632.825,13.37,SZ,,FBS,V2P1,07-26-2023"
$new.Range("J5").Value = "9,HP2,SZ-0526-3-3,sample,526,3,0526-3,3,1646.195,1632.825,13.37,SZ,,FBS,V2P1,07-26-2023"
$new.Range("K5").Value = "15,HP2,SZ-0602-1-3,sample,602,1,0602-1,3,1658.43,1646.84,11.59,SZ,,HPL,PLA,07-26-2023"
$new.Range("L5").Value = "15,HP2,SZ-0602-1-3,sample,602,1,0602-1,3,1658.43,1646.84,11.59,SZ,,HPL,PLA,07-26-2023"
$new.Range("C6").Value = "4,HP2,SZ-0526-2-1,sample,526,2,0526-2,1,1630.35,1618.18,12.17,SZ,,HPL,V2P1,07-26-2023"
$new.Range("D6").Value = "4,HP2,SZ-0526-2-1,sample,526,2,0526-2,1,1630.35,1618.18,12.17,SZ,,HPL,V2P1,07-26-2023"
$new.Range("E6").Value = "4,HP2,SZ-0526-2-1,sample,526,2,0526-2,1,1630.35,1618.18,12.17,SZ,,HPL,V2P1,07-26-2023"
$new.Range("F6").Value = "4,HP2,SZ-0526-2-1,sample,526,2,0526-2,1,1630.35,1618.18,12.17,SZ,,HPL,V2P1,07-26-2023"
$new.Range("G6").Value = "10,HP2,SZ-0526-4-1,sample,526,4,0526-4,1,1650.125,1636.775,13.35,SZ,,FBS,V2P1,07-26-2023"
$new.Range("H6").Value = "10,HP2,SZ-0526-4-1,sample,526,4,0526-4,1,1650.125,1636.775,13.35,SZ,,FBS,V2P1,07-26-2023"
$new.Range("I6").Value = "10,HP2,SZ-0526-4-1,sample,526,4,0526-4,1,1650.125,1636.775,13.35,SZ,,FBS,V2P1,07-26-2023"
$new.Range("J6").Value = "10,HP2,SZ-0526-4-1,sample,526,4,0526-4,1,1650.125,1636.775,13.35,SZ,,FBS,V2P1,07-26-2023"
$new.Range("K6").Value = "16,HP2,SZ-0602-2-1,sample,602,2,0602-2,1,1648.685,1638.37,10.315,SZ,,HPL,PLA,07-26-2023"
$new.Range("L6").Value = "16,HP2,SZ-0602-2-1,sample,602,2,0602-2,1,1648.685,1638.37,10.315,SZ,,HPL,PLA,07-26-2023"
$new.Range("C7").Value = "5,HP2,SZ-0526-2-2,sample,526,2,0526-2,2,1674.71,1661.275,13.435,SZ,,HPL,V2P1,07-26-2023"
$new.Range("D7").Value = "5,HP2,SZ-0526-2-2,sample,526,2,0526-2,2,1674.71,1661.275,13.435,SZ,,HPL,V2P1,07-26-2023"
$new.Range("E7").Value = "5,HP2,SZ-0526-2-2,sample,526,2,0526-2,2,1674.71,1661.275,13.435,SZ,,HPL,V2P1,07-26-2023"
$new.Range("F7").Value = "5,HP2,SZ-0526-2-2,sample,526,2,0526-2,2,1674.71,1661.275,13.435,SZ,,HPL,V2P1,07-26-2023"
$new.Range("G7").Value = "11,HP2,SZ-0526-4-2,sample,526,4,0526-4,2,1637.985,1624.165,13.82,SZ,,FBS,V2P1,07-26-2023"
$new.Range("H7").Value = "11,HP2,SZ-0526-4-2,sample,526,4,0526-4,2,1637.985,1624.165,13.82,SZ,,FBS,V2P1,07-26-2023"
$new.Range("I7").Value = "11,HP2,SZ-0526-4-2,sample,526,4,0526-4,2,1637.985,1624.165,13.82,SZ,,FBS,V2P1,07-26-2023"
$new.Range("J7").Value = "11,HP2,SZ-0526-4-2,sample,526,4,0526-4,2,1637.985,1624.165,13.82,SZ,,FBS,V2P1,07-26-2023"
$new.Range("K7").Value = "17,HP2,SZ-0602-2-2,sample,602,2,0602-2,2,1667.53,1656.865,10.665,SZ,,HPL,PLA,07-26-2023"
$new.Range("L7").Value = "17,HP2,SZ-0602-2-2,sample,602,2,0602-2,2,1667.53,1656.865,10.665,SZ,,HPL,PLA,07-26-2023"
$new.Range("C8").Value = "6,HP2,SZ-0526-2-3,sample,526,2,0526-2,3,1634.525,1621.96,12.565,SZ,,HPL,V2P1,07-26-2023"
$new.Range("D8").Value = "6,HP2,SZ-0526-2-3,sample,526,2,0526-2,3,1634.525,1621.96,12.565,SZ,,HPL,V2P1,07-26-2023"
$new.Range("E8").Value = "6,HP2,SZ-0526-2-3,sample,526,2,0526-2,3,1634.525,1621.96,12.565,SZ,,HPL,V2P1,07-26-2023"
$new.Range("F8").Value = "6,HP2,SZ-0526-2-3,sample,526,2,0526-2,3,1634.525,1621.96,12.565,SZ,,HPL,V2P1,07-26-2023"
$new.Range("G8").Value = "12,HP2,SZ-0526-4-3,sample,526,4,0526-4,3,1649.095,1635.38,13.715,SZ,,FBS,V2P1,07-26-2023"
$new.Range("H8").Value = "12,HP2,SZ-0526-4-3,sample,526,4,0526-4,3,1649.095,1635.38,13.715,SZ,,FBS,V2P1,07-26-2023"
$new.Range("I8").Value = "12,HP2,SZ-0526-4-3,sample,526,4,0526-4,3,1649.095,1635.38,13.715,SZ,,FBS,V2P1,07-26-2023"
$new.Range("J8").Value = "12,HP2,SZ-0526-4-3,sample,526,4,0526-4,3,1649.095,1635.38,13.715,SZ,,FBS,V2P1,07-26-2023"
$new.Range("K8").Value = "18,HP2,SZ-0602-2-3,sample,602,2,0602-2,3,1665.9,1655.535,10.365,SZ,,HPL,PLA,07-26-2023"
$new.Range("L8").Value = "18,HP2,SZ-0602-2-3,sample,602,2,0602-2,3,1665.9,1655.535,10.365,SZ,,HPL,PLA,07-26-2023"

# Match the saved selection/active cell shown in the source file.
[void]$new.Range("G9").Select()
